$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same event table and need the
# same update: a handful of "想去人数" (F column) counter bumps plus one new
# event row (row 14) appended at the bottom.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Bump the "want to go" counters in column F ---
    $ws.Range("F2").Value = 1513
    $ws.Range("F3").Value = 31
    $ws.Range("F4").Value = 977
    $ws.Range("F6").Value = 2388
    $ws.Range("F8").Value = 1466
    $ws.Range("F9").Value = 67
    $ws.Range("F10").Value = 163
    $ws.Range("F11").Value = 54
    $ws.Range("F12").Value = 408

    # --- Append the new event as row 14 ---
    # Column A uses the bold/bordered/centered style also used by the other
    # rows in column A; copy it from the row above so we reuse the existing
    # style record instead of minting a new (duplicate) one.
    $ws.Range("A13").Copy() | Out-Null
    $ws.Range("A14").PasteSpecial(-4122) | Out-Null
    $ws.Range("A14").Value = 13

    # Column B holds a literal text date like "2024.04.04". A plain
    # Range.Value assignment of that string gets auto-converted to a date
    # serial number, so instead build it as a text-returning formula and
    # immediately collapse it back to a plain value in place (copy + paste
    # special values-only) - this keeps the literal text "2024.04.04"
    # without leaving a formula behind or minting a new number-format style.
    $ws.Range("B14").Formula = '="2024.04.04"'
    $ws.Range("B14").Copy() | Out-Null
    $ws.Range("B14").PasteSpecial(-4163) | Out-Null

    $ws.Range("C14").Value = "赣州·赣次元·归来国风动漫节"
    $ws.Range("D14").Value = "客家大道568号文清外国语学校旁 赣州市文清外国语学校国际交流中心"
    $ws.Range("E14").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("F14").Value = 1
    $ws.Range("G14").Value = 40
    $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=82125"
    $ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202402/8RNepTak1709022774421.jpeg"
}
